# NEAP_Extent_DataOverlayMaster.xlsx edit
# "Change in overlay workflow for Marine input file"
#
# Row 3 (Marine input): rename dataset from "Natural Values Ecosystems" to
# "Marine" and point RawDataPath at the new intermediate Marine_Benthic.tif
# instead of the old NVE-Benthic raw tif.
#
# Rows 7/8 (ALUM inputs): rename "ALUM v8" entries to the year-specific
# "ALUM 2010" / "ALUM 2015" labels, drop the stale Crosswalk_comment on the
# 2015 row, and append a brand-new "ALUM 2020" row (9) wired up to the 2020-21
# ABARES land-use raster plus its resampled overlay-grid output, including a
# hyperlink on the ResampledRasterPath cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Marine input dataset swap -------------------------------------
$ws.Range("A3").Value = 'Marine'
$ws.Range("B3").Value = '\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\Marine_Benthic.tif'

# --- Row 7: ALUM 2010 (was "ALUM v8") --------------------------------------
$ws.Range("A7").Value = 'ALUM 2010'

# --- Row 8: ALUM 2015 (was "ALUM v8"); clear stale comment -----------------
$ws.Range("A8").Value = 'ALUM 2015'
$ws.Range("L8").Value = ""

# --- Row 9 (new): ALUM 2020 -------------------------------------------------
$ws.Range("A9").Value = 'ALUM 2020'
$ws.Range("B9").Value = '\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\inputs\raw\Land_use_of_Australia\ABARES_Land_use_of_Australia_2020_21_prerelease_20240712\ABARES_Land_use_of_Australia_2020_21_20240712\NLUM_v7p_ALUMV8_250m_2020_21_alb.tif'
$ws.Range("C9").Value = 'NA'
$ws.Range("D9").Value = 'Raster'
$ws.Range("E9").Value = 'epsg:3577'
$ws.Range("F9").Value = '250m'
$ws.Range("G9").Value = 'NN'
$ws.Range("H9").Value = 'NA'
$ws.Range("I9").Value = '\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\NEAP_NLUM_2020-21_EPSG3577_250m.tif'
$ws.Hyperlinks.Add($ws.Range("I9"), "file:///\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\overlay_grids\NEAP_NLUM_2020-21_EPSG3577_250m.tif")
$ws.Range("I9").Style = $ws.Range("I8").Style
$ws.Range("J9").Value = 'ALUM-IUCNGET'
$ws.Range("K9").Value = 'https://github.com/CSIRO-enviro-informatics/ecosystem-typology/raw/main/crosswalks/ALUM-IUCNGET/ALUM-IUCNGET.xlsx'

# --- Selection moved to A8 in the saved view -------------------------------
$ws.Range("A8").Select()
